# Extend the parallel-line contingency table from 14 columns (B:O) to 16
# columns (B:Q) and rebalance the I/K/M/O columns for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add the two new header values (matching O1's style) ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Rows 2-25: swap the I/K and M/O column values ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # column I
    $ws.Cells.Item($r, 11).Value = 1   # column K
    $ws.Cells.Item($r, 13).Value = 2   # column M
    $ws.Cells.Item($r, 15).Value = 1   # column O
    $ws.Cells.Item($r, 16).Value = 2   # column P (new)
    $ws.Cells.Item($r, 17).Value = 2   # column Q (new)
}
